# Regenerate save_data to use K instead of Strike# for stammen_craig 2022 sheet.
# Column G (header "K") gets recalculated values for most data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 0
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 2
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    36 = 1
    37 = 2
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
